$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: update status (E12) from "open" to "closed", add note text in G12,
# and increase the row height to fit the wrapped text.
$ws.Range("E12").Value = "closed"
$ws.Range("G12").Value = "显示是根据用户登陆的app store所在的地区显示的"
$ws.Rows.Item(12).RowHeight = 34

# Update the selection/scroll position on the sheet to match the saved view.
$ws.Range("E14").Select()
